$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.974.31"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.204.53"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.71"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.14"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.43"
$ws.Range("E10").Value = "  +6.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.59"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0779"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.113"
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.38"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.544.12"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.359.80"
$ws.Range("E16").Value = "  +7.03%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.76"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.734"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.884.00"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0885"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.20"
$ws.Range("E21").Value = "  -2.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.75"
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.46"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.00"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.48"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.06"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.09"
$ws.Range("E29").Value = "  -4.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.28"
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.46"
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.80"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.94"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.03"
$ws.Range("E35").Value = "  +5.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0711"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0997"
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.51"
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.073.71"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.75"
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.20"
$ws.Range("E44").Value = "  +7.03%  "
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.86"
$ws.Range("E46").Value = "  -1.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.75"
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("E48").Value = "  -10.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.418.63"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.12"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.46"
$ws.Range("E51").Value = "  +1.03%  "
